# Add team record (Wins / Losses / Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1=Wins, AE1=Losses, AF1=Ties ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell (AC1) onto the new header cells so they match the rest of
# row 1's style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-42): team record for every player row ---
$wins = 71
$losses = 91
$ties = 0

for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
